# Update "想去人数" (want-to-go count) figures in the F column for the
# 展览 (sheet1) and 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 435
$ws1.Range("F4").Value = 1526
$ws1.Range("F5").Value = 8818
$ws1.Range("F9").Value = 308
$ws1.Range("F10").Value = 163
$ws1.Range("F13").Value = 3682
$ws1.Range("F17").Value = 3139
$ws1.Range("F22").Value = 2469
$ws1.Range("F23").Value = 84

# Sheet "全部类型" (All categories)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 435
$ws4.Range("F4").Value = 1526
$ws4.Range("F5").Value = 8818
$ws4.Range("F9").Value = 308
$ws4.Range("F10").Value = 163
$ws4.Range("F13").Value = 3682
$ws4.Range("F17").Value = 3139
$ws4.Range("F22").Value = 2469
$ws4.Range("F24").Value = 84
